$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("w01_100.5_optimization_results")

$data = @(
    @(1, 628, 0.7238711426878679, -1.30414857083974, -0.8927336325794407, 1.648687535547951),
    @(2, 638, 1.302221338218278, -0.1479775544684723, -0.4747336934124086, 2.979241101441685),
    @(3, 632, -0.6761042507329533, -1.077394989064373, -0.9990562842355675, 1.023872590304876),
    @(4, 638, 1.258547615003693, -0.01, -1, 1.001),
    @(5, 638, 1.094852883182382, -1.043444337938924, -0.4220136303950404, 2.940611963135439),
    @(6, 614, 1.065258100229737, -1.691024642156352, -0.9482972362502432, 2.744659194566054),
    @(7, 615, 0.1235064299036601, -1.85896629303264, -0.5763369310095678, 2.399892543246757),
    @(8, 621, 0.7719757117229418, -0.1001311950360371, -0.6959169693034615, 2.887677251564717),
    @(9, 622, -0.02567476635498211, -1.600340218309853, -0.66577904615033, 2.525892064905803),
    @(10, 622, -1.130787719431383, -0.4796247492573122, -0.4207786728043932, 1.015247883023973)
)

$startRow = 12
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
    $ws.Cells.Item($row, 5).Value = $values[4]
    $ws.Cells.Item($row, 6).Value = $values[5]
}
